# Add a new worksheet named "1787" at the end of the workbook, cloning the
# layout of the existing "695" sheet (same labels/structure) but using
# "Swerve" instead of "Tank" for the robot-type row (A10).

$wb = $excel.ActiveWorkbook

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "1787"

# Row 1 header cells
$newSheet.Range("C1").Value = "Ranks"
$newSheet.Range("E1").Value = "Broke?"
$newSheet.Range("F1").Value = "Match #"
$newSheet.Range("G1").Value = "Scouter"
$newSheet.Range("H1").Value = "Comment"

# Column A labels (rows 2-10)
$newSheet.Range("A2").Value = "Total Points"
$newSheet.Range("A3").Value = "Tele Points"
$newSheet.Range("A4").Value = "Auto Points"
$newSheet.Range("A5").Value = "Coral Points"
$newSheet.Range("A6").Value = "Algae Points"
$newSheet.Range("A7").Value = "Rice Score"
$newSheet.Range("A8").Value = "--------------------------------------------------------------------------------------------------"
$newSheet.Range("A9").Value = "Qualitative"
$newSheet.Range("A10").Value = "Swerve"

# Row 39 header cells
$newSheet.Range("A39").Value = "Match"
$newSheet.Range("B39").Value = "Auto"
$newSheet.Range("C39").Value = "Tele"
$newSheet.Range("D39").Value = "Person"
